# Actualización automática hashcode lun oct  5 00:53:10 CEST 2020
# Update the hashcode values (column B) for the rows listed below.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hashcode.csv")

$updates = @{
    9   = "cca34025ec355a3510b22aed4feb0be0"
    17  = "8f3320683070b351b8a46525d9398033"
    126 = "e3975995e93715cac89992b3c6abf56d"
    136 = "6872be526233ba8c3b6fb36d445b4acf"
    183 = "4ddd5dd1470070e4d517b35f8bdeda84"
    200 = "34b67f98ec9c468b8c9de264ff0af310"
    227 = "85989131ddf30ca9b10b0f700aa4b08d"
    228 = "a86a111bfa73f6fa5114926803991516"
    232 = "ad18caeb5133af7054afc04f4b1190a8"
    284 = "9f863b1993aed11c7455afa38fb7e995"
    305 = "152ea32c5063267575bb6dbcb40a4a90"
    467 = "8dea7e76a3223ee22747bd732db4307c"
    486 = "bd9c8ff20b6463a37026ed37b3a7fc9f"
    510 = "45294000247a48924ff590fa40686167"
    562 = "f077fe0b721ce58d7095be780bc1df1d"
    584 = "83864ca16281fdb70c702b6bbeb12d8a"
    712 = "f350b655bac9cf43cc7698210390daed"
    715 = "c067500944de81d06d9ed5fd4b75776f"
    727 = "956718bd0775f4546661f78b0b233033"
    842 = "ec20e7fccf7b4f694af9f796acc044de"
    847 = "5ad6b3d876277d316a39f9048e40b49e"
    874 = "9a9f3f6ebea23c1f1bb3896d197e98f3"
}

foreach ($row in $updates.Keys) {
    $ws.Range("B$row").Value = $updates[$row]
}
